# Auto-generated update: appends new tracker rows 151-166 (matches commit "Actualizacion automatica del tracker")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template empty cells (G150:H150) are copied down for each new row so that the
# newly appended G/H ("resultado"/"profit") cells exist as genuine empty cells,
# matching the still-empty inlineStr cells used for not-yet-settled picks.
$emptyTemplate = $ws.Range("G150:H150")

$ws.Cells.Item(151, 1).Value = 14684189
$ws.Cells.Item(151, 2).NumberFormat = "@"
$ws.Cells.Item(151, 2).Value = "2025-09-16"
$ws.Cells.Item(151, 3).Value = "Diana Shnaider"
$ws.Cells.Item(151, 4).Value = "Catherine McNally"
$ws.Cells.Item(151, 5).Value = "Gana Diana Shnaider"
$ws.Cells.Item(151, 6).Value = 1.62
$emptyTemplate.Copy($ws.Range("G151:H151"))

$ws.Cells.Item(152, 1).Value = 14680555
$ws.Cells.Item(152, 2).NumberFormat = "@"
$ws.Cells.Item(152, 2).Value = "2025-09-15"
$ws.Cells.Item(152, 3).Value = "Francesco Passaro"
$ws.Cells.Item(152, 4).Value = "Dusan Lajovic"
$ws.Cells.Item(152, 5).Value = "Gana Dusan Lajovic"
$ws.Cells.Item(152, 6).Value = 1.8
$emptyTemplate.Copy($ws.Range("G152:H152"))

$ws.Cells.Item(153, 1).Value = 14680548
$ws.Cells.Item(153, 2).NumberFormat = "@"
$ws.Cells.Item(153, 2).Value = "2025-09-16"
$ws.Cells.Item(153, 3).Value = "Sebastian Sorger"
$ws.Cells.Item(153, 4).Value = "Jerome Kym"
$ws.Cells.Item(153, 5).Value = "Gana Sebastian Sorger"
$ws.Cells.Item(153, 6).Value = 6
$emptyTemplate.Copy($ws.Range("G153:H153"))

$ws.Cells.Item(154, 1).Value = 14686286
$ws.Cells.Item(154, 2).NumberFormat = "@"
$ws.Cells.Item(154, 2).Value = "2025-09-16"
$ws.Cells.Item(154, 3).Value = "Cezar Stefan Bentzel"
$ws.Cells.Item(154, 4).Value = "Gabi Adrian Boitan"
$ws.Cells.Item(154, 5).Value = "Gana Cezar Stefan Bentzel"
$ws.Cells.Item(154, 6).Value = 15
$emptyTemplate.Copy($ws.Range("G154:H154"))

$ws.Cells.Item(155, 1).Value = 14685777
$ws.Cells.Item(155, 2).NumberFormat = "@"
$ws.Cells.Item(155, 2).Value = "2025-09-15"
$ws.Cells.Item(155, 3).Value = "Michael Vrbensky"
$ws.Cells.Item(155, 4).Value = "Tadeas Paroulek"
$ws.Cells.Item(155, 5).Value = "Gana Michael Vrbensky"
$ws.Cells.Item(155, 6).Value = 1.57
$emptyTemplate.Copy($ws.Range("G155:H155"))

$ws.Cells.Item(156, 1).Value = 14686832
$ws.Cells.Item(156, 2).NumberFormat = "@"
$ws.Cells.Item(156, 2).Value = "2025-09-15"
$ws.Cells.Item(156, 3).Value = "Trey Hilderbrand"
$ws.Cells.Item(156, 4).Value = "Aryan Shah"
$ws.Cells.Item(156, 5).Value = "Gana Trey Hilderbrand"
$ws.Cells.Item(156, 6).Value = 3.5
$emptyTemplate.Copy($ws.Range("G156:H156"))

$ws.Cells.Item(157, 1).Value = 14687249
$ws.Cells.Item(157, 2).NumberFormat = "@"
$ws.Cells.Item(157, 2).Value = "2025-09-15"
$ws.Cells.Item(157, 3).Value = "Antoine Ghibaudo"
$ws.Cells.Item(157, 4).Value = "Erik Arutiunian"
$ws.Cells.Item(157, 5).Value = "Gana Erik Arutiunian"
$ws.Cells.Item(157, 6).Value = 2.25
$emptyTemplate.Copy($ws.Range("G157:H157"))

$ws.Cells.Item(158, 1).Value = 14687251
$ws.Cells.Item(158, 2).NumberFormat = "@"
$ws.Cells.Item(158, 2).Value = "2025-09-15"
$ws.Cells.Item(158, 3).Value = "Samir Banerjee"
$ws.Cells.Item(158, 4).Value = "Nikita Samuel Filin"
$ws.Cells.Item(158, 5).Value = "Gana Nikita Samuel Filin"
$ws.Cells.Item(158, 6).Value = 8
$emptyTemplate.Copy($ws.Range("G158:H158"))

$ws.Cells.Item(159, 1).Value = 14687911
$ws.Cells.Item(159, 2).NumberFormat = "@"
$ws.Cells.Item(159, 2).Value = "2025-09-15"
$ws.Cells.Item(159, 3).Value = "Daniel Milavsky"
$ws.Cells.Item(159, 4).Value = "Sebastian Dominko"
$ws.Cells.Item(159, 5).Value = "Gana Sebastian Dominko"
$ws.Cells.Item(159, 6).Value = 2.25
$emptyTemplate.Copy($ws.Range("G159:H159"))

$ws.Cells.Item(160, 1).Value = 14686096
$ws.Cells.Item(160, 2).NumberFormat = "@"
$ws.Cells.Item(160, 2).Value = "2025-09-15"
$ws.Cells.Item(160, 3).Value = "Gianluca Cadenasso"
$ws.Cells.Item(160, 4).Value = "Alejandro Mateo Berge Nourescu"
$ws.Cells.Item(160, 5).Value = "Gana Alejandro Mateo Berge Nourescu"
$ws.Cells.Item(160, 6).Value = 7
$emptyTemplate.Copy($ws.Range("G160:H160"))

$ws.Cells.Item(161, 1).Value = 14686089
$ws.Cells.Item(161, 2).NumberFormat = "@"
$ws.Cells.Item(161, 2).Value = "2025-09-15"
$ws.Cells.Item(161, 3).Value = "Gabriele Pennaforti"
$ws.Cells.Item(161, 4).Value = "Imanol Lopez Morillo"
$ws.Cells.Item(161, 5).Value = "Gana Imanol Lopez Morillo"
$ws.Cells.Item(161, 6).Value = 2
$emptyTemplate.Copy($ws.Range("G161:H161"))

$ws.Cells.Item(162, 1).Value = 14686091
$ws.Cells.Item(162, 2).NumberFormat = "@"
$ws.Cells.Item(162, 2).Value = "2025-09-15"
$ws.Cells.Item(162, 3).Value = "Tommaso Compagnucci"
$ws.Cells.Item(162, 4).Value = "Stefan Palosi"
$ws.Cells.Item(162, 5).Value = "Gana Stefan Palosi"
$ws.Cells.Item(162, 6).Value = 1.57
$emptyTemplate.Copy($ws.Range("G162:H162"))

$ws.Cells.Item(163, 1).Value = 14685801
$ws.Cells.Item(163, 2).NumberFormat = "@"
$ws.Cells.Item(163, 2).Value = "2025-09-15"
$ws.Cells.Item(163, 3).Value = "Hernan Casanova"
$ws.Cells.Item(163, 4).Value = "Alejo Lorenzo Lingua Lavallen"
$ws.Cells.Item(163, 5).Value = "Gana Alejo Lorenzo Lingua Lavallen"
$ws.Cells.Item(163, 6).Value = 4.33
$emptyTemplate.Copy($ws.Range("G163:H163"))

$ws.Cells.Item(164, 1).Value = 14686701
$ws.Cells.Item(164, 2).NumberFormat = "@"
$ws.Cells.Item(164, 2).Value = "2025-09-15"
$ws.Cells.Item(164, 3).Value = "Facundo Bagnis"
$ws.Cells.Item(164, 4).Value = "Leonardo Aboian"
$ws.Cells.Item(164, 5).Value = "Gana Leonardo Aboian"
$ws.Cells.Item(164, 6).Value = 6.5
$emptyTemplate.Copy($ws.Range("G164:H164"))

$ws.Cells.Item(165, 1).Value = 14686933
$ws.Cells.Item(165, 2).NumberFormat = "@"
$ws.Cells.Item(165, 2).Value = "2025-09-15"
$ws.Cells.Item(165, 3).Value = "Ignacio Monzon"
$ws.Cells.Item(165, 4).Value = "Carlos Maria Zarate"
$ws.Cells.Item(165, 5).Value = "Gana Carlos Maria Zarate"
$ws.Cells.Item(165, 6).Value = 1.83
$emptyTemplate.Copy($ws.Range("G165:H165"))

$ws.Cells.Item(166, 1).Value = 14685452
$ws.Cells.Item(166, 2).NumberFormat = "@"
$ws.Cells.Item(166, 2).Value = "2025-09-15"
$ws.Cells.Item(166, 3).Value = "Vitalia Diatchenko"
$ws.Cells.Item(166, 4).Value = "Carol Young-suh Lee"
$ws.Cells.Item(166, 5).Value = "Gana Carol Young-suh Lee"
$ws.Cells.Item(166, 6).Value = 1.83
$emptyTemplate.Copy($ws.Range("G166:H166"))

